$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$p = $d.Paragraphs.Item(1)
$r = $p.Range
$r.End = $r.End - 1  # exclude the paragraph mark
$r.Text = "2024-01-22 Monday"

# Update the arithmetic expressions in the table, cell by cell.
# We set each cell Range.Text directly (rather than Find/Replace)
# because several expressions repeat verbatim elsewhere in the table,
# and a document-wide Find/Replace could hit the wrong occurrence.
$tbl = $d.Tables.Item(1)

$cell = $tbl.Rows.Item(1).Cells.Item(1)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "51-48="
$cell = $tbl.Rows.Item(1).Cells.Item(2)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "33-20="
$cell = $tbl.Rows.Item(1).Cells.Item(3)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "31+26="
$cell = $tbl.Rows.Item(1).Cells.Item(4)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "33+50="
$cell = $tbl.Rows.Item(1).Cells.Item(5)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "88-64="
$cell = $tbl.Rows.Item(2).Cells.Item(1)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "24+26="
$cell = $tbl.Rows.Item(2).Cells.Item(2)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "74-22="
$cell = $tbl.Rows.Item(2).Cells.Item(3)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "27+67="
$cell = $tbl.Rows.Item(2).Cells.Item(4)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "76-72="
$cell = $tbl.Rows.Item(2).Cells.Item(5)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "37+29="
$cell = $tbl.Rows.Item(3).Cells.Item(1)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "3+46="
$cell = $tbl.Rows.Item(3).Cells.Item(2)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "87-0="
$cell = $tbl.Rows.Item(3).Cells.Item(3)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "86-43="
$cell = $tbl.Rows.Item(3).Cells.Item(4)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "93-30="
$cell = $tbl.Rows.Item(3).Cells.Item(5)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "6-3="
$cell = $tbl.Rows.Item(4).Cells.Item(1)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "51+38="
$cell = $tbl.Rows.Item(4).Cells.Item(2)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "57-2="
$cell = $tbl.Rows.Item(4).Cells.Item(3)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "58+36="
$cell = $tbl.Rows.Item(4).Cells.Item(4)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "20+22="
$cell = $tbl.Rows.Item(4).Cells.Item(5)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "88-39="
$cell = $tbl.Rows.Item(5).Cells.Item(1)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "52+35="
$cell = $tbl.Rows.Item(5).Cells.Item(2)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "55-24="
$cell = $tbl.Rows.Item(5).Cells.Item(3)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "97-70="
$cell = $tbl.Rows.Item(5).Cells.Item(4)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "78+9="
$cell = $tbl.Rows.Item(5).Cells.Item(5)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "65-26="
$cell = $tbl.Rows.Item(6).Cells.Item(1)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "92-75="
$cell = $tbl.Rows.Item(6).Cells.Item(2)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "98-84="
$cell = $tbl.Rows.Item(6).Cells.Item(3)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "0+66="
$cell = $tbl.Rows.Item(6).Cells.Item(4)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "31-8="
$cell = $tbl.Rows.Item(6).Cells.Item(5)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "37+19="
$cell = $tbl.Rows.Item(7).Cells.Item(1)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "61-60="
$cell = $tbl.Rows.Item(7).Cells.Item(2)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "26-0="
$cell = $tbl.Rows.Item(7).Cells.Item(3)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "31-7="
$cell = $tbl.Rows.Item(7).Cells.Item(4)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "19+64="
$cell = $tbl.Rows.Item(7).Cells.Item(5)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "28+2="
$cell = $tbl.Rows.Item(8).Cells.Item(1)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "13+78="
$cell = $tbl.Rows.Item(8).Cells.Item(2)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "64-47="
$cell = $tbl.Rows.Item(8).Cells.Item(3)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "91-19="
$cell = $tbl.Rows.Item(8).Cells.Item(4)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "6+5="
$cell = $tbl.Rows.Item(8).Cells.Item(5)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "34+44="
$cell = $tbl.Rows.Item(9).Cells.Item(1)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "94-94="
$cell = $tbl.Rows.Item(9).Cells.Item(2)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "5+0="
$cell = $tbl.Rows.Item(9).Cells.Item(3)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "82-68="
$cell = $tbl.Rows.Item(9).Cells.Item(4)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "81-29="
$cell = $tbl.Rows.Item(9).Cells.Item(5)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "4+28="
$cell = $tbl.Rows.Item(10).Cells.Item(1)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "29+1="
$cell = $tbl.Rows.Item(10).Cells.Item(2)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "96-52="
$cell = $tbl.Rows.Item(10).Cells.Item(3)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "51+25="
$cell = $tbl.Rows.Item(10).Cells.Item(4)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "50+5="
$cell = $tbl.Rows.Item(10).Cells.Item(5)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "4+27="
$cell = $tbl.Rows.Item(11).Cells.Item(1)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "44-20="
$cell = $tbl.Rows.Item(11).Cells.Item(2)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "87-73="
$cell = $tbl.Rows.Item(11).Cells.Item(3)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "71-27="
$cell = $tbl.Rows.Item(11).Cells.Item(4)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "71-36="
$cell = $tbl.Rows.Item(11).Cells.Item(5)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "94-71="
$cell = $tbl.Rows.Item(12).Cells.Item(1)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "52+23="
$cell = $tbl.Rows.Item(12).Cells.Item(2)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "74-73="
$cell = $tbl.Rows.Item(12).Cells.Item(3)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "97-88="
$cell = $tbl.Rows.Item(12).Cells.Item(4)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "76-1="
$cell = $tbl.Rows.Item(12).Cells.Item(5)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "74-12="
$cell = $tbl.Rows.Item(13).Cells.Item(1)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "65+20="
$cell = $tbl.Rows.Item(13).Cells.Item(2)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "7+58="
$cell = $tbl.Rows.Item(13).Cells.Item(3)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "22+30="
$cell = $tbl.Rows.Item(13).Cells.Item(4)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "78-40="
$cell = $tbl.Rows.Item(13).Cells.Item(5)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "69-7="
$cell = $tbl.Rows.Item(14).Cells.Item(1)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "65-10="
$cell = $tbl.Rows.Item(14).Cells.Item(2)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "69+0="
$cell = $tbl.Rows.Item(14).Cells.Item(3)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "15+24="
$cell = $tbl.Rows.Item(14).Cells.Item(4)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "32+66="
$cell = $tbl.Rows.Item(14).Cells.Item(5)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "5+33="
$cell = $tbl.Rows.Item(15).Cells.Item(1)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "66-2="
$cell = $tbl.Rows.Item(15).Cells.Item(2)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "25-24="
$cell = $tbl.Rows.Item(15).Cells.Item(3)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "89-28="
$cell = $tbl.Rows.Item(15).Cells.Item(4)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "59-9="
$cell = $tbl.Rows.Item(15).Cells.Item(5)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "79-40="
$cell = $tbl.Rows.Item(16).Cells.Item(1)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "5+64="
$cell = $tbl.Rows.Item(16).Cells.Item(2)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "33-3="
$cell = $tbl.Rows.Item(16).Cells.Item(3)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "46+39="
$cell = $tbl.Rows.Item(16).Cells.Item(4)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "63-19="
$cell = $tbl.Rows.Item(16).Cells.Item(5)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "38-20="
$cell = $tbl.Rows.Item(17).Cells.Item(1)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "82+15="
$cell = $tbl.Rows.Item(17).Cells.Item(2)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "10+81="
$cell = $tbl.Rows.Item(17).Cells.Item(3)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "65-25="
$cell = $tbl.Rows.Item(17).Cells.Item(4)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "82-59="
$cell = $tbl.Rows.Item(17).Cells.Item(5)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "8+22="
$cell = $tbl.Rows.Item(18).Cells.Item(1)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "78+14="
$cell = $tbl.Rows.Item(18).Cells.Item(2)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "0+68="
$cell = $tbl.Rows.Item(18).Cells.Item(3)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "3+30="
$cell = $tbl.Rows.Item(18).Cells.Item(4)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "76-6="
$cell = $tbl.Rows.Item(18).Cells.Item(5)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "20-4="
$cell = $tbl.Rows.Item(19).Cells.Item(1)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "20+0="
$cell = $tbl.Rows.Item(19).Cells.Item(2)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "86-24="
$cell = $tbl.Rows.Item(19).Cells.Item(3)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "41+54="
$cell = $tbl.Rows.Item(19).Cells.Item(4)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "77-1="
$cell = $tbl.Rows.Item(19).Cells.Item(5)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "32+33="
$cell = $tbl.Rows.Item(20).Cells.Item(1)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "55-39="
$cell = $tbl.Rows.Item(20).Cells.Item(2)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "1+7="
$cell = $tbl.Rows.Item(20).Cells.Item(3)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "89-89="
$cell = $tbl.Rows.Item(20).Cells.Item(4)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "67-62="
$cell = $tbl.Rows.Item(20).Cells.Item(5)
$cr = $cell.Range
$cr.End = $cr.End - 1  # exclude the cell-end mark
$cr.Text = "1+53="
